# Update cryptos list with latest prices and percentage changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.052.28'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '4.013.18'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '525.68'
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.89'
$ws.Range("E6").Value = '  +1.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.697'
$ws.Range("E7").Value = '  +12.55%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.746'
$ws.Range("E9").Value = '  +1.76%  '
$ws.Range("E10").Value = '  -1.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '50.94'
$ws.Range("E11").Value = '  +9.64%  '
$ws.Range("E12").Value = '  -2.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.71'
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("D14").Value = '4.651.58'
$ws.Range("E14").Value = '  -0.84%  '
$ws.Range("D15").Value = '4.021.42'
$ws.Range("E15").Value = '  -0.92%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.99'
$ws.Range("E16").Value = '  -1.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.54'
$ws.Range("E17").Value = '  -2.78%  '
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.19'
$ws.Range("E19").Value = '  -2.01%  '
$ws.Range("D20").Value = '71.878.88'
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '429.78'
$ws.Range("E21").Value = '  -1.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '97.21'
$ws.Range("E22").Value = '  +2.33%  '
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.17'
$ws.Range("E24").Value = '  +2.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.28'
$ws.Range("E25").Value = '  -0.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.09'
$ws.Range("E26").Value = '  -8.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.70'
$ws.Range("E27").Value = '  -4.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.86'
$ws.Range("E28").Value = '  +1.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.68'
$ws.Range("E29").Value = '  +20.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.57'
$ws.Range("E30").Value = '  -0.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.42'
$ws.Range("E31").Value = '  +6.68%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.131'
$ws.Range("E32").Value = '  +1.56%  '
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '13.39'
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '679.22'
$ws.Range("E34").Value = '  -3.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '47.42'
$ws.Range("E35").Value = '  +16.95%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '65.37'
$ws.Range("E36").Value = '  -3.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.446'
$ws.Range("E37").Value = '  +1.14%  '
$ws.Range("E38").Value = '  -1.77%  '
$ws.Range("D39").Value = '0.0₃0822'
$ws.Range("E39").Value = '  -8.73%  '
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.39'
$ws.Range("E40").Value = '  +8.33%  '
$ws.Range("B41").Value = 'ThetaToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.38'
$ws.Range("E41").Value = '  -8.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0488'
$ws.Range("E44").Value = '  +1.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.13'
$ws.Range("E45").Value = '  +12.31%  '
$ws.Range("E46").Value = '  +3.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.65'
$ws.Range("E47").Value = '  -4.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.41'
$ws.Range("E48").Value = '  -3.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.99'
$ws.Range("E49").Value = '  -4.85%  '
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000268'
$ws.Range("E50").Value = '  -3.46%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '144.25'
$ws.Range("E51").Value = '  +1.84%  '
